$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Finalise UI." / " (2 hours)" were two separate runs - merge them
#    into a single run reading "Finalise UI. (2 hours)" (only the
#    first occurrence in the document needs this - a later, unrelated
#    paragraph already has the merged text).
# ---------------------------------------------------------------------
$null = $d.Content.Find.Execute("Finalise UI. (2 hours)", $true, $false, $false, $false, $false, `
                                 $true, 1, $false, "Finalise UI. (2 hours)", 1)

# ---------------------------------------------------------------------
# 2) Drop the old "_GoBack" bookmark that sat after "Help with
#    application testing." - Word will re-create it at the new last
#    edit location below.
# ---------------------------------------------------------------------
$oldGoBack = $d.Bookmarks.Item("_GoBack")
$oldGoBack.Delete()

# ---------------------------------------------------------------------
# 3) Make sure a "Hyperlink" character style exists (Word auto-creates
#    this the first time a hyperlink is inserted).
# ---------------------------------------------------------------------
$hlStyle = $d.Styles.Add("Hyperlink", 2)
$hlStyle.BaseStyle = "DefaultParagraphFont"
$hlStyle.Priority = 99
$hlStyle.UnhideWhenUsed = $true
$hlStyle.Font.Color = 16711680
$hlStyle.Font.Underline = 1

# ---------------------------------------------------------------------
# 4) Append two new paragraphs at the very end of the document:
#      - an (empty) paragraph that will hold the bookmarkStart
#      - "Project GitHub " followed by a hyperlink to the repo
# ---------------------------------------------------------------------
$endOfDoc = $d.Content.End
$tail = $d.Range($endOfDoc - 1, $endOfDoc - 1)
$tail.InsertParagraphAfter()
$bmPara = $d.Paragraphs.Last
$bmPara.Style = "Normal"
$bmPara.SpaceAfter = 0
$bmPara.LineSpacingRule = 0

# temporary placeholder character so the upcoming bookmark range has a
# well defined start strictly inside this paragraph (it is stripped
# out again once the bookmark has been created)
$bmParaRange = $bmPara.Range
$bmParaRange.Collapse(1)
$bmParaRange.Text = "X"

$afterBmPara = $bmPara.Range
$afterBmPara.Collapse(0)
$afterBmPara.InsertParagraphAfter()
$linkPara = $d.Paragraphs.Last
$linkPara.Style = "Normal"
$linkPara.SpaceAfter = 0
$linkPara.LineSpacingRule = 0

$linkParaRange = $linkPara.Range
$linkParaRange.Collapse(0)
$linkUrl = "https://github.com/crouchbindset/professionalskills"
$linkParaRange.Text = "Project GitHub " + $linkUrl

$urlStart = $linkPara.Range.Start + ("Project GitHub ").Length
$urlEnd = $linkPara.Range.End - 1
$urlRange = $d.Range($urlStart, $urlEnd)
$null = $d.Hyperlinks.Add($urlRange, $linkUrl)

# ---------------------------------------------------------------------
# 5) Re-create the "_GoBack" bookmark spanning from the placeholder
#    paragraph through to the end of the hyperlink paragraph, then
#    remove the placeholder character - this leaves bookmarkStart
#    alone in the empty paragraph and bookmarkEnd right after the
#    hyperlink, matching Word's own "last edit" bookkeeping.
# ---------------------------------------------------------------------
$newBmStart = $bmPara.Range.Start
$newBmEnd = $linkPara.Range.End - 1
$newBmRange = $d.Range($newBmStart, $newBmEnd)
$null = $d.Bookmarks.Add("_GoBack", $newBmRange)

$placeholder = $d.Range($newBmStart, $newBmStart + 1)
$placeholder.Text = ""

Write-Output "done"
